$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.363.61"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.712.60"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'224.70"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "0.5291"
$ws.Range("E6").Value = "  -1.39%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D11").Value = "0.07668"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "4.514"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").Value = "1.948.39"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "1.712.81"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "0.5836"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "0.0₅8229"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "'68.20"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "27.366.55"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "223.38"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'4.630"
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("E22").Value = "  -2.40%  "
$ws.Range("D23").Value = "6.016"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "144.97"
$ws.Range("E25").Value = "  -2.16%  "
$ws.Range("D26").Value = "1.689"
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("D28").Value = "7.243"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").Value = "16.31"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").Value = "0.05363"
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "3.485"
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("D33").Value = "3.437"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").Value = "2.873"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").Value = "0.9511"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").Value = "0.5858"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").Value = "0.01636"
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").Value = "1.095.74"
$ws.Range("E40").Value = "  +3.77%  "
$ws.Range("D41").Value = "5.793"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "0.8408"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").Value = "101.11"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").Value = "1.855.03"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").Value = "0.4537"
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("D49").Value = "1.003"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").Value = "8.115"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "0.05239"
$ws.Range("E51").Value = "  -0.23%  "
